# Added environmental variable support for creds
# Append 5 new rows (41-45) to the "stock_data" sheet, each a copy of the
# last existing data row (row 40) but dated one day later ("Oct 06, 2022").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

# Values mirror row 40 verbatim, except column A which gets the new date.
$rowValues = @(
    "Oct 06, 2022", "NYSE",
    3342, 1339, 1828, 175, 18, 101,
    294788518, 478710051, 789929201,
    1.03,
    3890,
    1485040244, 2080024643,
    "NASDAQ",
    3665857390,
    4813, 1696, 2811, 306, 42, 216,
    1.26,
    21021,
    1359398333, 2842883736, 4257416449
)

for ($r = 41; $r -le 45; $r++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $rowValues[$i]
    }
}
